{"js": "// \"Increased instant access limit\"\n//\n// The Instant Access CPUh limit was raised from 7000 to 20000. The figure\n// is mentioned twice in the document body:\n//   1. \"You can apply for a maximum of 7000 CPUh.\"\n//   2. \"Successful instant access projects receive up to 7000 CPUh ...\"\n// Both need \"7000\" changed to \"20000\"; no other visible text changes.\n\nconst body = context.document.body;\n\n// 1. Replace every occurrence of \"7000\" with \"20000\".\nconst hits = body.search(\"7000\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < hits.items.length; i++) {\n  hits.items[i].insertText(\"20000\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2. Word tracks the location of the most recent edit with a hidden\n//    \"_GoBack\" bookmark. Before this edit it sat next to an unrelated\n//    paragraph; after typing the replacement it naturally lands right\n//    after \"...receive up to 20\" (i.e. where the user's cursor ended up\n//    replacing \"7\" with \"20\" in the second sentence). Move it to match.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst marker = context.document.body.search(\"receive up to 20\", { matchCase: true });\nmarker.load(\"items\");\nawait context.sync();\n\nif (marker.items.length > 0) {\n  const afterMatch = marker.items[0].getRange(Word.RangeLocation.after);\n  afterMatch.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"Increased instant access limit\"\n#\n# The Instant Access CPUh limit was raised from 7000 to 20000. The figure\n# is mentioned twice in the document body:\n#   1. \"You can apply for a maximum of 7000 CPUh.\"\n#   2. \"Successful instant access projects receive up to 7000 CPUh ...\"\n# Both need \"7000\" changed to \"20000\"; no other visible text changes.\n\n$d = $word.ActiveDocument\n\n# 1. Replace every occurrence of \"7000\" with \"20000\".\n$find = $d.Content.Find\n$find.Text = \"7000\"\n$find.Replacement.Text = \"20000\"\n$find.Execute(\n    $find.Text,              # FindText\n    $false,                  # MatchCase\n    $false,                  # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                  # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                        # Replace (wdReplaceAll)\n)\n\n# 2. Word tracks the location of the most recent edit with a hidden\n#    \"_GoBack\" bookmark. Before this edit it sat next to an unrelated\n#    paragraph; after typing the replacement it naturally lands right\n#    after \"...receive up to 20\" (i.e. where the user's cursor ended up\n#    replacing \"7\" with \"20\" in the second sentence). Move it to match.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$markerFind = $d.Content.Find\n$markerFind.Text = \"receive up to 20\"\n$markerFind.Execute()\nif ($markerFind.Found) {\n    $markerEnd = $markerFind.Parent.End\n    $bmRange = $d.Range($markerEnd, $markerEnd)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
